$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and Column E (Volume 1h) updates for the refreshed cryptos snapshot.
# Numeric-looking price strings must be forced back to text (matching the original
# inline-string cell type) so Excel does not silently coerce them into numbers and
# drop significant trailing zeros / change formatting.

$ws.Range('D2').Value = '42.554.68'
$ws.Range('E2').Value = '  -1.15%  '
$ws.Range('D3').Value = '2.344.14'
$ws.Range('E3').Value = '  -1.72%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.94'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -3.80%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.03'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.632'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -2.24%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.614'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -6.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.27'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0924'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.61'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.72%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.995'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -5.33%  '
$ws.Range('E14').Value = '  +0.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.94'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -7.41%  '
$ws.Range('D16').Value = '2.699.50'
$ws.Range('E16').Value = '  -1.62%  '
$ws.Range('D17').Value = '2.331.17'
$ws.Range('E17').Value = '  -1.88%  '
$ws.Range('D18').Value = '42.484.09'
$ws.Range('E18').Value = '  -1.16%  '
$ws.Range('E19').Value = '  -3.49%  '
$ws.Range('E20').Value = '  -2.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '75.70'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.69%  '
$ws.Range('E22').Value = '  +0.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '255.30'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -6.23%  '
$ws.Range('E24').Value = '  -4.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.42'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.38'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -3.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.65'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.34%  '
$ws.Range('E29').Value = '  +1.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '172.64'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '36.95'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0890'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.03'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.53%  '
$ws.Range('E34').Value = '  -9.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.123'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +15.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.132'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.98%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.63'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -6.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0361'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.41%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.92'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -6.15%  '
$ws.Range('E40').Value = '  -5.61%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.239'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.25%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.46'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -7.66%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '70.43'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.14%  '
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.03'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.59%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '110.96'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -8.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.18'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.73%  '
$ws.Range('E48').Value = '  -1.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '85.07'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -6.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '74.76'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.29'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.80%  '
